# Add a new forecast column BB to the sheet, one column to the right of
# the previous last column BA, mirroring the existing forecast layout:
#  - Row 1 gets a new period-end date header (BB1)
#  - Rows 3-18 carry forward the same (already-converged) forecast value
#    that is currently in column BA for that row
#  - Rows 19-21 get a new, distinct forecast value for the new period
#  - Rows 2 and 22 are left untouched (no BB value)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header for column BB, following the same style as BA1 (s="1")
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("BB1").Value = 45986

# Rows 3-18: copy forward the converged BA value into the new BB column
$carryForwardRows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18)
foreach ($r in $carryForwardRows) {
    $baAddr = "BA" + $r
    $bbAddr = "BB" + $r
    $ws.Range($bbAddr).Value = $ws.Range($baAddr).Value2
}

# Rows 19-21: new forecast values for the newly added period
$ws.Range("BB19").Value = 2.622852459381209
$ws.Range("BB20").Value = 1.946625946175717
$ws.Range("BB21").Value = 2.030414939481551

$excel.CutCopyMode = $false
